$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Insert a new "Headline" worksheet right after "CTA" (before "RTE"),
#    mirroring the other authoring-component test sheets (Hero / CTA).
# ---------------------------------------------------------------------
$cta = $wb.Worksheets.Item("CTA")
$headline = $wb.Worksheets.Add($null, $cta)
$headline.Name = "Headline"

# Header row
$headline.Range("A1").Value = "TestName"
$headline.Range("B1").Value = "headlineText"
$headline.Range("C1").Value = "headlineLink"
$headline.Range("D1").Value = "headlineLinkOption"

# Row 2: plain-text headline sample
$headline.Range("A2").Value = "PlainText"
$headline.Range("B2").Value = "Sample Test"
$headline.Range("D2").Value = "Existing window/tab"

# Row 3: linked headline sample
$headline.Range("A3").Value = "Link"
$headline.Range("B3").Value = "Sample Test"
$headline.Range("C3").Value = "/content/pathology-education/language-masters/en/testing"
$headline.Range("D3").Value = "New tab"

# Approximate the authored column widths for the new sheet
$headline.Columns.Item(2).ColumnWidth = 11.54296875
$headline.Columns.Item(3).ColumnWidth = 23.26953125
$headline.Columns.Item(4).ColumnWidth = 17.1796875

# Leave the selection on D2:D3, matching the authored sheet state
$headline.Range("D2:D3").Select()

# ---------------------------------------------------------------------
# 2. Update the "TestCases" sheet: the Hero smoke-test row now points at
#    the new Headline component, and the old stray CTA count row is
#    removed.
# ---------------------------------------------------------------------
$testCases = $wb.Worksheets.Item("TestCases")
$testCases.Rows.Item(3).Delete()
$testCases.Range("A2").Value = "Headline"
$testCases.Range("B2").Select()

# TestCases becomes the active/selected sheet again
$testCases.Select()
